# Weekly update: a new price observation is inserted as row 51 (pushing the
# existing rows 51-127 down to 52-128), extending the data range from
# A1:R127 to A1:R128.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 51; Excel shifts rows 51..127 down to 52..128
# and grows the used range accordingly (A1:R127 -> A1:R128).
$ws.Rows("51").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A51").Value = 10
$ws.Range("B51").Value = "Vega Modelo de Temuco"
$ws.Range("C51").Value = "La Araucanía"
$ws.Range("D51").Value = 44579
$ws.Range("E51").Value = 9
$ws.Range("F51").Value = 100114007
$ws.Range("G51").Value = "Jengibre"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 50
$ws.Range("K51").Value = 20000
$ws.Range("L51").Value = 20000
$ws.Range("M51").Value = 20000
$ws.Range("N51").Value = "$/caja 13 kilos"
$ws.Range("O51").Value = "Perú"
$ws.Range("P51").Value = 1538
$ws.Range("Q51").Value = 13
$ws.Range("R51").Value = "Hortaliza"
